# This script applies the diff: updates the date line and the
# 3-digit x 1-digit multiplication problems/answers throughout the table.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-01-01 Wednesday"; New = "2025-01-02 Thursday" },
    @{ Old = "522×9=4698"; New = "284×9=2556" },
    @{ Old = "182×9=1638"; New = "771×8=6168" },
    @{ Old = "251×3=753";  New = "694×3=2082" },
    @{ Old = "912×6=5472"; New = "339×6=2034" },
    @{ Old = "363×4=1452"; New = "590×9=5310" },
    @{ Old = "988×8=7904"; New = "700×6=4200" },
    @{ Old = "613×2=1226"; New = "794×4=3176" },
    @{ Old = "312×2=624";  New = "495×9=4455" },
    @{ Old = "257×7=1799"; New = "453×6=2718" },
    @{ Old = "990×2=1980"; New = "710×9=6390" },
    @{ Old = "291×8=2328"; New = "677×4=2708" },
    @{ Old = "842×9=7578"; New = "219×5=1095" },
    @{ Old = "826×2=1652"; New = "284×2=568" },
    @{ Old = "930×5=4650"; New = "846×3=2538" },
    @{ Old = "395×7=2765"; New = "404×2=808" },
    @{ Old = "352×8=2816"; New = "638×7=4466" },
    @{ Old = "291×6=1746"; New = "768×8=6144" },
    @{ Old = "702×7=4914"; New = "583×5=2915" },
    @{ Old = "657×2=1314"; New = "167×9=1503" },
    @{ Old = "601×7=4207"; New = "153×2=306" },
    @{ Old = "568×8=4544"; New = "228×5=1140" },
    @{ Old = "859×6=5154"; New = "583×8=4664" },
    @{ Old = "835×9=7515"; New = "674×3=2022" },
    @{ Old = "323×3=969";  New = "758×9=6822" },
    @{ Old = "735×7=5145"; New = "522×3=1566" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
